$wb = $excel.ActiveWorkbook

# Sheet ALC, row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3515.5833
$ws.Range("I40").Value = 2774.75
$ws.Range("J40").Value = 4997.25
$ws.Range("K40").Value = 2774.75
$ws.Range("L40").Value = 4997.25
$ws.Range("M40").Value = -2599.75
$ws.Range("N40").Value = -5347.25

# Sheet ALC, row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3484
$ws.Range("I76").Value = 2902.3333
$ws.Range("J76").Value = 4880
$ws.Range("K76").Value = 2902.3333
$ws.Range("L76").Value = 4880
$ws.Range("M76").Value = -2587.3333
$ws.Range("N76").Value = -5510

# Sheet ALC, row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3484
$ws.Range("I79").Value = 2902.3333
$ws.Range("J79").Value = 4880
$ws.Range("K79").Value = 2902.3333
$ws.Range("L79").Value = 4880
$ws.Range("M79").Value = -1810.3333
$ws.Range("N79").Value = -7064

# Sheet ALC, row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2257.6428
$ws.Range("I106").Value = 2123.6155
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 2123.6155
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -1492.6155
$ws.Range("N106").Value = -5262

# Sheet ALC, row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6249
$ws.Range("I113").Value = 1874
$ws.Range("J113").Value = 7999
$ws.Range("K113").Value = 1874
$ws.Range("L113").Value = 7999
$ws.Range("M113").Value = 1380
$ws.Range("N113").Value = -14507

# Sheet ALC, row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3123251.2
$ws.Range("I132").Value = 3192457
$ws.Range("J132").Value = 9000
$ws.Range("K132").Value = 9577371
$ws.Range("L132").Value = 27000
$ws.Range("M132").Value = -9574841
$ws.Range("N132").Value = -32060

# Sheet ALC, row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 15497.583
$ws.Range("I137").Value = 20202.646
$ws.Range("J137").Value = 4071
$ws.Range("K137").Value = 60607.938
$ws.Range("L137").Value = 12213
$ws.Range("M137").Value = -58057.938
$ws.Range("N137").Value = -17313

# Sheet ALC, row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2818.0244
$ws.Range("I138").Value = 1736.2727
$ws.Range("J138").Value = 4070.5789
$ws.Range("K138").Value = 5208.8181
$ws.Range("L138").Value = 12211.7367
$ws.Range("M138").Value = -68.81810000000041
$ws.Range("N138").Value = -22491.7367

# Sheet ALC, row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2863.375
$ws.Range("I141").Value = 2770.4783
$ws.Range("K141").Value = 8311.4349
$ws.Range("M141").Value = -3131.4349

# Sheet ARM, row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18129.309
$ws.Range("I32").Value = 18642.969
$ws.Range("K32").Value = 18642.969
$ws.Range("M32").Value = -18355.969

# Sheet ARM, row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 291989.06
$ws.Range("I74").Value = 317119.47
$ws.Range("J74").Value = 53250
$ws.Range("K74").Value = 317119.47
$ws.Range("L74").Value = 53250
$ws.Range("M74").Value = -316245.47
$ws.Range("N74").Value = -54998

# Sheet ARM, row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 291989.06
$ws.Range("I77").Value = 317119.47
$ws.Range("J77").Value = 53250
$ws.Range("K77").Value = 1585597.35
$ws.Range("L77").Value = 266250
$ws.Range("M77").Value = -1581229.35
$ws.Range("N77").Value = -274986

# Sheet ARM, row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1998.2046
$ws.Range("I122").Value = 1994.55
$ws.Range("J122").Value = 2034.75
$ws.Range("K122").Value = 5983.65
$ws.Range("L122").Value = 6104.25
$ws.Range("M122").Value = -3533.65
$ws.Range("N122").Value = -11004.25

# Sheet ARM, row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1345.36
$ws.Range("I132").Value = 1244.9565
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 3734.8695
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -1204.8695
$ws.Range("N132").Value = -12560

# Sheet BSM, row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2398.182
$ws.Range("I99").Value = 2138
$ws.Range("K99").Value = 2138
$ws.Range("M99").Value = -640

# Sheet CRP, row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5002785
$ws.Range("I31").Value = 8334640.5
$ws.Range("J31").Value = 5001.25
$ws.Range("K31").Value = 8334640.5
$ws.Range("L31").Value = 5001.25
$ws.Range("M31").Value = -8334345.5
$ws.Range("N31").Value = -5591.25

# Sheet CRP, row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5002785
$ws.Range("I34").Value = 8334640.5
$ws.Range("J34").Value = 5001.25
$ws.Range("K34").Value = 8334640.5
$ws.Range("L34").Value = 5001.25
$ws.Range("M34").Value = -8334438.5
$ws.Range("N34").Value = -5405.25

# Sheet CRP, row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 829.1111
$ws.Range("I94").Value = 952.1667
$ws.Range("J94").Value = 583
$ws.Range("K94").Value = 952.1667
$ws.Range("L94").Value = 583
$ws.Range("M94").Value = -501.1667
$ws.Range("N94").Value = -1485

# Sheet CUL, row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2778.1667
$ws.Range("I5").Value = 3193.8
$ws.Range("J5").Value = 700
$ws.Range("K5").Value = 9581.400000000001
$ws.Range("L5").Value = 2100
$ws.Range("M5").Value = -9469.400000000001
$ws.Range("N5").Value = -2324

# Sheet CUL, row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1061.4166
$ws.Range("I98").Value = 416.75
$ws.Range("J98").Value = 1383.75
$ws.Range("K98").Value = 1250.25
$ws.Range("L98").Value = 4151.25
$ws.Range("M98").Value = 247.75
$ws.Range("N98").Value = -7147.25

# Sheet CUL, row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1886.8572
$ws.Range("I132").Value = 2655.1667
$ws.Range("J132").Value = 1579.5333
$ws.Range("K132").Value = 23896.5003
$ws.Range("L132").Value = 14215.7997
$ws.Range("M132").Value = -21366.5003
$ws.Range("N132").Value = -19275.7997

# Sheet CUL, row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2778.1667
$ws.Range("I135").Value = 3193.8
$ws.Range("J135").Value = 700
$ws.Range("K135").Value = 28744.2
$ws.Range("L135").Value = 6300
$ws.Range("M135").Value = -26209.2
$ws.Range("N135").Value = -11370

# Sheet CUL, row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5488.8
$ws.Range("I137").Value = 5309.4287
$ws.Range("J137").Value = 8000
$ws.Range("K137").Value = 15928.2861
$ws.Range("L137").Value = 24000
$ws.Range("M137").Value = -10828.2861
$ws.Range("N137").Value = -34200

# Sheet GSM, row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 17721.742
$ws.Range("I102").Value = 18788.793
$ws.Range("J102").Value = 2249.5
$ws.Range("K102").Value = 18788.793
$ws.Range("L102").Value = 2249.5
$ws.Range("M102").Value = -17166.793
$ws.Range("N102").Value = -5493.5

# Sheet GSM, row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3103.25
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 3837.6667
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 3837.6667
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -8177.6667

# Sheet GSM, row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1859.875
$ws.Range("I132").Value = 1889.5
$ws.Range("J132").Value = 1534
$ws.Range("K132").Value = 5668.5
$ws.Range("L132").Value = 4602
$ws.Range("M132").Value = -3138.5
$ws.Range("N132").Value = -9662

# Sheet LTW, row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3170.1714
$ws.Range("I122").Value = 2998.4333
$ws.Range("J122").Value = 4200.6
$ws.Range("K122").Value = 8995.2999
$ws.Range("L122").Value = 12601.8
$ws.Range("M122").Value = -6545.2999
$ws.Range("N122").Value = -17501.8

# Sheet LTW, row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2544.75
$ws.Range("I132").Value = 2349.1765
$ws.Range("J132").Value = 4539.6
$ws.Range("K132").Value = 7047.529500000001
$ws.Range("L132").Value = 13618.8
$ws.Range("M132").Value = -4517.529500000001
$ws.Range("N132").Value = -18678.8

# Sheet LTW, row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3438.0476
$ws.Range("I136").Value = 3483.25
$ws.Range("J136").Value = 3166.8333
$ws.Range("K136").Value = 10449.75
$ws.Range("L136").Value = 9500.499899999999
$ws.Range("M136").Value = -7899.75
$ws.Range("N136").Value = -14600.4999

# Sheet WVR, row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1114.5
$ws.Range("I113").Value = 905.3077
$ws.Range("J113").Value = 1416.6666
$ws.Range("K113").Value = 2715.9231
$ws.Range("L113").Value = 4249.9998
$ws.Range("M113").Value = -545.9231
$ws.Range("N113").Value = -8589.9998

# Sheet WVR, row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18136.125
$ws.Range("I132").Value = 20372.096
$ws.Range("J132").Value = 2484.3333
$ws.Range("K132").Value = 61116.288
$ws.Range("L132").Value = 7452.999899999999
$ws.Range("M132").Value = -58586.288
$ws.Range("N132").Value = -12512.9999

# Sheet WVR, row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 14311.135
$ws.Range("I136").Value = 16103.934
$ws.Range("J136").Value = 2786
$ws.Range("K136").Value = 48311.802
$ws.Range("L136").Value = 8358
$ws.Range("M136").Value = -45761.802
$ws.Range("N136").Value = -13458
